# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Matches the commit: "Created functions to get season record" — every
# player row gets the team's season record appended in columns AD:AF.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new labels in AD1:AF1, styled like the rest of the
# header row (bold, centered, bordered) by copying A1's format over.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows (2-47): season record is constant for every player on the team.
$wins = 85
$losses = 77
$ties = 1

$lastRow = 47
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
